# Set up the run mode for the test suite.
#
# Insert a new worksheet named "test_suite" as the very first sheet of the
# workbook. Worksheets.Add() (with no arguments) inserts the new sheet
# immediately before the currently active sheet, which - since the active
# sheet is the first one - places it at position 1, pushing the existing
# "AddCustomerTest" / "OpenAccountTest" sheets down.
$wb = $excel.ActiveWorkbook
$testSuite = $wb.Worksheets.Add()
$testSuite.Name = "test_suite"

# Header row.
$testSuite.Range("A1").Value = "TCID"
$testSuite.Range("B1").Value = "Runmode"

# Fill in the TCID column for the two pre-existing test cases first ...
$testSuite.Range("A2").Value = "AddCustomerTest"
$testSuite.Range("A3").Value = "OpenAccountTest"

# ... then insert a brand-new row above them for a third test case so it
# ends up listed first in the suite.
$testSuite.Rows("2:2").Insert()
$testSuite.Range("A2").Value = "BankManagerLoginTest"

# Now fill in the Runmode column.
$testSuite.Range("B2").Value = "Y"
$testSuite.Range("B3").Value = "N"
$testSuite.Range("B4").Value = "N"

# Size column A to fit its contents (best effort - mirrors Excel's own
# "best fit" autosize that fires after typing the TCID column) and leave
# the selection on B3, matching where the author last edited.
$testSuite.Columns("A:A").ColumnWidth = 21
$testSuite.Range("B3").Select() | Out-Null
